$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (sending cluster order changed: ECs, FAPs, M1, M2, sCs)
# and all numeric columns recomputed per Dr Hou's advice.
$data = @(
    @("ECs",  "Pdgfb", "Art1", "sCs", 3, 1,              49.195614,          147.586842,         0.8724243622956617,   0.8724243622956616,   3, 1, 0.1748566666666667, 0.52457, 1, 1, 8.602181078659999,  77.41962970793999, 0.8724243622956617,   0.8724243622956616),
    @("FAPs", "Pdgfb", "Art1", "sCs", 1, 0.3333333333333333, 0.05229733333333333, 0.156892,       0.0009274295810821058, 0.0009274295810821058, 3, 1, 0.1748566666666667, 0.52457, 1, 1, 0.009144537382222222, 0.08230083644,      0.0009274295810821058, 0.0009274295810821058),
    @("M1",   "Pdgfb", "Art1", "sCs", 3, 1,              1.896931,           5.690793,           0.03363976345521111,  0.0336397634552111,   3, 1, 0.1748566666666667, 0.52457, 1, 1, 0.3316910315566667, 2.98521928401,     0.03363976345521111,  0.0336397634552111),
    @("M2",   "Pdgfb", "Art1", "sCs", 3, 1,              3.047190333333333, 9.141570999999999,  0.05403821331210212,  0.05403821331210212,  3, 1, 0.1748566666666667, 0.52457, 1, 1, 0.5328215443855555, 4.79539389947,     0.05403821331210212,  0.05403821331210212),
    @("sCs",  "Pdgfb", "Art1", "sCs", 3, 1,              2.197513666666667, 6.592541,            0.03897023135594298,  0.03897023135594298,  3, 1, 0.1748566666666667, 0.52457, 1, 1, 0.3842499147077778, 3.45824923237,     0.03897023135594298,  0.03897023135594298)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
